# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet to the latest scraped snapshot. Matches the GitHub Actions commit
# "Updated cryptos list on Mon May 15 10:56:17 UTC 2023 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '27.668.94'; E = '  -0.54%  ' },
    @{ Row = 3; D = '1.846.31'; E = '  -1.11%  ' },
    @{ Row = 4; D = '1.010'; E = '  -2.90%  ' },
    @{ Row = 5; D = '319.12'; E = '  -1.79%  ' },
    @{ Row = 6; D = '1.010'; E = '  -2.56%  ' },
    @{ Row = 7; D = '0.4310'; E = '  -2.65%  ' },
    @{ Row = 8; D = '0.3747'; E = '  -1.49%  ' },
    @{ Row = 9; D = '0.07350'; E = '  -1.83%  ' },
    @{ Row = 10; D = '0.8806'; E = '  -0.69%  ' },
    @{ Row = 11; D = '21.58'; E = '  -0.96%  ' },
    @{ Row = 12; D = '1.840.28'; E = '  -1.73%  ' },
    @{ Row = 13; D = '6.735'; E = '  -0.50%  ' },
    @{ Row = 14; D = '5.453'; E = '  -2.07%  ' },
    @{ Row = 15; D = '0.07120'; E = '  -1.61%  ' },
    @{ Row = 16; D = '87.70'; E = '  +4.50%  ' },
    @{ Row = 17; D = '1.014'; E = '  -2.76%  ' },
    @{ Row = 18; D = '0.000008989'; E = '  -2.03%  ' },
    @{ Row = 19; D = $null; E = '  -2.63%  ' },
    @{ Row = 20; D = '15.49'; E = '  -0.56%  ' },
    @{ Row = 21; D = '27.677.16'; E = '  -0.57%  ' },
    @{ Row = 22; D = '5.254'; E = '  -1.45%  ' },
    @{ Row = 23; D = '11.17'; E = '  -1.96%  ' },
    @{ Row = 24; D = '2.081.23'; E = '  -1.45%  ' },
    @{ Row = 25; D = '2.036'; E = '  +1.93%  ' },
    @{ Row = 26; D = '155.58'; E = '  -2.04%  ' },
    @{ Row = 27; D = '18.58'; E = '  -1.70%  ' },
    @{ Row = 28; D = '2.143'; E = '  +7.68%  ' },
    @{ Row = 29; D = '5.386'; E = '  +0.70%  ' },
    @{ Row = 30; D = '120.45'; E = '  +2.17%  ' },
    @{ Row = 31; D = '0.08923'; E = '  -1.79%  ' },
    @{ Row = 32; D = '1.231'; E = '  +1.10%  ' },
    @{ Row = 33; D = '0.7806'; E = '  +0.00%  ' },
    @{ Row = 34; D = '4.560'; E = '  -0.57%  ' },
    @{ Row = 35; D = '2.909'; E = '  -7.09%  ' },
    @{ Row = 36; D = $null; E = '  -2.70%  ' },
    @{ Row = 37; D = '1.139'; E = '  -1.56%  ' },
    @{ Row = 38; D = '0.05332'; E = $null },
    @{ Row = 39; D = '0.01970'; E = '  -1.47%  ' },
    @{ Row = 40; D = '7.237'; E = '  +4.18%  ' },
    @{ Row = 41; D = '2.861'; E = '  -0.54%  ' },
    @{ Row = 42; D = '0.5163'; E = '  -0.99%  ' },
    @{ Row = 43; D = '0.1679'; E = '  -1.10%  ' },
    @{ Row = 44; D = '8.937'; E = '  +2.69%  ' },
    @{ Row = 45; D = '110.28'; E = '  +0.49%  ' },
    @{ Row = 46; D = '10.62'; E = '  -1.44%  ' },
    @{ Row = 47; D = '0.4732'; E = '  +0.14%  ' },
    @{ Row = 48; D = '0.06496'; E = '  +0.38%  ' },
    @{ Row = 49; D = '1.699'; E = '  -1.82%  ' },
    @{ Row = 50; D = '1.011'; E = '  -2.80%  ' },
    @{ Row = 51; D = '1.889'; E = '  -1.22%  ' }
)

# The Price column (D) holds values such as "27.660.83", "1.011" or
# "0.07344" that are stored as plain text (locale-style thousands
# separators, leading zeros, etc.). Excel's COM layer auto-detects plain
# numeric-looking strings assigned via .Value and silently converts them
# to a Double (losing precision/formatting), so the whole column is
# force-formatted as Text before writing, then restored to the sheet's
# normal (General) style so no stray formatting is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

$priceRange.Style = "Normal"
